# Added save functionality and updated documentation
# Replace the numeric indicator matrix (rows 2-5, cols A-D) with text labels:
# a cell keeps an empty string where the old value was 0, and gets a text
# label where the old value was 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = ""
$ws.Range("B2").Value = "krawedz"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = "inna_krawedz"

$ws.Range("A3").Value = ""
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = "krawedz"
$ws.Range("D3").Value = ""

$ws.Range("A4").Value = ""
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "krawedz"

$ws.Range("A5").Value = ""
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = ""
